$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44251
$ws.Range("J2").Value2 = 120
$ws.Range("L2").Value2 = 5000
$ws.Range("M2").Value2 = 5000
$ws.Range("O2").Value2 = 'Región Metropolitana'
$ws.Range("P2").Value2 = 312

# Row 3
$ws.Range("D3").Value2 = 44208

# Row 4
$ws.Range("D4").Value2 = 44188
$ws.Range("J4").Value2 = 210
$ws.Range("L4").Value2 = 6000
$ws.Range("M4").Value2 = 5500
$ws.Range("O4").Value2 = 'Provincia de Quillota'
$ws.Range("P4").Value2 = 344

# Row 5
$ws.Range("D5").Value2 = 44187
$ws.Range("J5").Value2 = 160

# Row 6
$ws.Range("D6").Value2 = 44232

# Row 7
$ws.Range("D7").Value2 = 44215
$ws.Range("J7").Value2 = 250

# Row 8
$ws.Range("D8").Value2 = 44210
$ws.Range("J8").Value2 = 340

# Row 9
$ws.Range("D9").Value2 = 44230
$ws.Range("J9").Value2 = 250

# Row 10
$ws.Range("D10").Value2 = 44292
$ws.Range("J10").Value2 = 90
$ws.Range("K10").Value2 = 6000
$ws.Range("M10").Value2 = 6000
$ws.Range("O10").Value2 = 'Región Metropolitana'
$ws.Range("P10").Value2 = 375

# Row 11
$ws.Range("D11").Value2 = 44204
$ws.Range("J11").Value2 = 430

# Row 12
$ws.Range("D12").Value2 = 44231
$ws.Range("J12").Value2 = 250
$ws.Range("K12").Value2 = 5000
$ws.Range("M12").Value2 = 5500
$ws.Range("O12").Value2 = 'Provincia de Quillota'
$ws.Range("P12").Value2 = 344

# Row 13
$ws.Range("D13").Value2 = 44186
$ws.Range("J13").Value2 = 160

# Row 14
$ws.Range("D14").Value2 = 44189
